# Auto-generated edit script applying the diff to Asura_Profits sheets
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 5035.25
$ws.Range("I11").Value = 5035.25
$ws.Range("K11").Value = 5035.25
$ws.Range("M11").Value = -4895.25
$ws.Range("H31").Value = 12279.25
$ws.Range("I31").Value = 12279.25
$ws.Range("K31").Value = 36837.75
$ws.Range("M31").Value = -36607.75
$ws.Range("H34").Value = 17709.777
$ws.Range("I34").Value = 17709.777
$ws.Range("K34").Value = 17709.777
$ws.Range("M34").Value = -17506.777
$ws.Range("H36").Value = 17709.777
$ws.Range("I36").Value = 17709.777
$ws.Range("K36").Value = 17709.777
$ws.Range("M36").Value = -16994.777
$ws.Range("H88").Value = 1904.7778
$ws.Range("J88").Value = 2513.3333
$ws.Range("L88").Value = 2513.3333
$ws.Range("N88").Value = -3325.3333
$ws.Range("H91").Value = 1904.7778
$ws.Range("J91").Value = 2513.3333
$ws.Range("L91").Value = 2513.3333
$ws.Range("N91").Value = -5321.3333
$ws.Range("H106").Value = 8500
$ws.Range("I106").Value = 10000
$ws.Range("J106").Value = 7000
$ws.Range("K106").Value = 10000
$ws.Range("L106").Value = 7000
$ws.Range("M106").Value = -9369
$ws.Range("N106").Value = -8262
$ws.Range("H125").Value = 3769.353
$ws.Range("J125").Value = 3790.7693
$ws.Range("L125").Value = 34116.9237
$ws.Range("N125").Value = -39036.9237
$ws.Range("H138").Value = 3146.2942
$ws.Range("I138").Value = 1512.0333
$ws.Range("J138").Value = 4436.5
$ws.Range("K138").Value = 4536.0999
$ws.Range("L138").Value = 13309.5
$ws.Range("M138").Value = 603.9000999999998
$ws.Range("N138").Value = -23589.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("H24").Value = 28885
$ws.Range("J24").Value = 28885
$ws.Range("L24").Value = 28885
$ws.Range("N24").Value = -29633
$ws.Range("H100").Value = 28885
$ws.Range("J100").Value = 28885
$ws.Range("L100").Value = 28885
$ws.Range("N100").Value = -31049

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 5249
$ws.Range("I5").Value = 5249
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 5249
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -5136
$ws.Range("N5").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 9999.5
$ws.Range("J4").Value = 9999.5
$ws.Range("L4").Value = 9999.5
$ws.Range("N4").Value = -10223.5
$ws.Range("H31").Value = 1968.2642
$ws.Range("I31").Value = 1224.7727
$ws.Range("J31").Value = 5603.1113
$ws.Range("K31").Value = 1224.7727
$ws.Range("L31").Value = 5603.1113
$ws.Range("M31").Value = -929.7727
$ws.Range("N31").Value = -6193.1113
$ws.Range("H34").Value = 1968.2642
$ws.Range("I34").Value = 1224.7727
$ws.Range("J34").Value = 5603.1113
$ws.Range("K34").Value = 1224.7727
$ws.Range("L34").Value = 5603.1113
$ws.Range("M34").Value = -1022.7727
$ws.Range("N34").Value = -6007.1113
$ws.Range("H129").Value = 49999
$ws.Range("J129").Value = 49999
$ws.Range("L129").Value = 49999
$ws.Range("N129").Value = -59999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 1836.2
$ws.Range("I13").Value = 1795.25
$ws.Range("J13").Value = 2000
$ws.Range("K13").Value = 5385.75
$ws.Range("L13").Value = 6000
$ws.Range("M13").Value = -5217.75
$ws.Range("N13").Value = -6336
$ws.Range("H58").Value = 2316.75
$ws.Range("J58").Value = 2316.75
$ws.Range("L58").Value = 6950.25
$ws.Range("N58").Value = -7206.25
$ws.Range("H64").Value = 3100
$ws.Range("J64").Value = 3230
$ws.Range("L64").Value = 9690
$ws.Range("N64").Value = -10230
$ws.Range("H67").Value = 3100
$ws.Range("J67").Value = 3230
$ws.Range("L67").Value = 9690
$ws.Range("N67").Value = -11562
$ws.Range("H70").Value = 16242.4
$ws.Range("J70").Value = 10000
$ws.Range("L70").Value = 30000
$ws.Range("N70").Value = -30630
$ws.Range("H73").Value = 16242.4
$ws.Range("J73").Value = 10000
$ws.Range("L73").Value = 30000
$ws.Range("N73").Value = -32184
$ws.Range("H93").Value = 1300
$ws.Range("J93").Value = 1500
$ws.Range("L93").Value = 4500
$ws.Range("N93").Value = -8244
$ws.Range("H107").Value = 1057.909
$ws.Range("I107").Value = 1300
$ws.Range("J107").Value = 1004.1111
$ws.Range("K107").Value = 3900
$ws.Range("L107").Value = 3012.3333
$ws.Range("M107").Value = -1980
$ws.Range("N107").Value = -6852.3333
$ws.Range("H108").Value = 3468.7
$ws.Range("I108").Value = 981.25
$ws.Range("J108").Value = 4090.5625
$ws.Range("K108").Value = 2943.75
$ws.Range("L108").Value = 12271.6875
$ws.Range("M108").Value = -63.75
$ws.Range("N108").Value = -18031.6875
$ws.Range("H109").Value = 100027
$ws.Range("I109").Value = 100027
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 300081
$ws.Range("L109").Value = 0
$ws.Range("M109").Value = -299041
$ws.Range("N109").ClearContents()
$ws.Range("H119").Value = 6943
$ws.Range("I119").Value = 552.6667
$ws.Range("K119").Value = 1658.0001
$ws.Range("M119").Value = 3179.9999
$ws.Range("H122").Value = 450.8
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H123").Value = 8330
$ws.Range("H131").Value = 1299.2273
$ws.Range("I131").Value = 747.1429000000001
$ws.Range("J131").Value = 1403.6757
$ws.Range("K131").Value = 2241.4287
$ws.Range("L131").Value = 4211.0271
$ws.Range("M131").Value = 2798.5713
$ws.Range("N131").Value = -14291.0271

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H59").Value = 50000
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 35000
$ws.Range("J42").Value = 30000
$ws.Range("L42").Value = 30000
$ws.Range("N42").Value = -31126
$ws.Range("H49").Value = 35000
$ws.Range("J49").Value = 30000
$ws.Range("L49").Value = 30000
$ws.Range("N49").Value = -30294

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 40000
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 40000
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 40000
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -40826
$ws.Range("H42").Value = 38340.668
$ws.Range("I42").Value = 38348
$ws.Range("J42").Value = 38333.332
$ws.Range("K42").Value = 38348
$ws.Range("L42").Value = 38333.332
$ws.Range("M42").Value = -37970
$ws.Range("N42").Value = -39089.332
$ws.Range("H43").Value = 23333.334
$ws.Range("H51").Value = 14792.8
$ws.Range("I51").Value = 5991
$ws.Range("K51").Value = 5991
$ws.Range("M51").Value = -5481
